$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Candidatures")

# Row 4 (Charlie C) is duplicated into the new row 5, only the name differs
# ("Denise D" instead of "Charlie C") -- everything else (postes, cycle,
# cours donnes, nobels, discipline, cote Z, choix) stays identical, along
# with the row's formatting.
$ws.Range("A4:H4").Copy() | Out-Null
$ws.Range("A5:H5").PasteSpecial() | Out-Null

$ws.Range("A5").Value = "Denise D"

# Restore the view: scroll back to the top-left and leave the selection on
# the cell right below the newly added row.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A6").Select() | Out-Null
